$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.683.49'
$ws.Range('E2').Value = '  +0.62%  '

$ws.Range('D3').Value = '1.820.18'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '''228.53'
$ws.Range('E5').Value = '  +0.65%  '

$ws.Range('D6').Value = '''0.584'
$ws.Range('E6').Value = '  +5.27%  '

$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').Value = '''34.81'
$ws.Range('E8').Value = '  +7.24%  '

$ws.Range('E9').Value = '  +1.35%  '

$ws.Range('E10').Value = '  +0.93%  '

$ws.Range('E11').Value = '  +0.35%  '

$ws.Range('D12').Value = '2.080.87'
$ws.Range('E12').Value = '  +1.17%  '

$ws.Range('D13').Value = '''11.45'
$ws.Range('E13').Value = '  +3.21%  '

$ws.Range('D14').Value = '1.824.61'
$ws.Range('E14').Value = '  +1.62%  '

$ws.Range('D16').Value = '34.657.90'

$ws.Range('D17').Value = '''4.34'
$ws.Range('E17').Value = '  +1.94%  '

$ws.Range('D18').Value = '''69.23'
$ws.Range('E18').Value = '  +0.77%  '

$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('E19').Value = '  +0.10%  '

$ws.Range('D20').Value = '''246.29'
$ws.Range('E20').Value = '  -0.26%  '

$ws.Range('D21').Value = '''11.57'
$ws.Range('E21').Value = '  +3.50%  '

$ws.Range('E22').Value = '  +0.14%  '

$ws.Range('E23').Value = '  +0.26%  '

$ws.Range('D24').Value = '''173.30'
$ws.Range('E24').Value = '  +5.94%  '

$ws.Range('E25').Value = '  +1.27%  '

$ws.Range('E26').Value = '  +3.02%  '

$ws.Range('E27').Value = '  +1.47%  '

$ws.Range('D28').Value = '''0.120'
$ws.Range('E28').Value = '  +3.07%  '

$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('E30').Value = '  +1.87%  '

$ws.Range('D31').Value = '''0.0530'
$ws.Range('E31').Value = '  +1.42%  '

$ws.Range('D32').Value = '''3.85'
$ws.Range('E32').Value = '  +1.69%  '

$ws.Range('E33').Value = '  +1.06%  '

$ws.Range('E34').Value = '  +0.94%  '

$ws.Range('D35').Value = '''2.61'
$ws.Range('E35').Value = '  -0.11%  '

$ws.Range('D36').Value = '1.407.05'

$ws.Range('D37').Value = '''0.682'
$ws.Range('E37').Value = '  +1.78%  '

$ws.Range('D38').Value = '''1.06'
$ws.Range('E38').Value = '  +0.25%  '

$ws.Range('D39').Value = '''0.0191'
$ws.Range('E39').Value = '  +0.06%  '

$ws.Range('D40').Value = '''84.05'
$ws.Range('E40').Value = '  -0.11%  '

$ws.Range('E41').Value = '  +4.98%  '

$ws.Range('D42').Value = '''0.952'
$ws.Range('E42').Value = '  +1.98%  '

$ws.Range('E43').Value = '  -0.13%  '

$ws.Range('D44').Value = '''13.75'
$ws.Range('E44').Value = '  -0.43%  '

$ws.Range('E45').Value = '  +2.88%  '

$ws.Range('E46').Value = '  -1.57%  '

$ws.Range('E47').Value = '  -0.36%  '

$ws.Range('D48').Value = '1.982.50'
$ws.Range('E48').Value = '  +1.47%  '

$ws.Range('D49').Value = '''105.34'
$ws.Range('E49').Value = '  -0.40%  '

$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.21%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0129'
$ws.Range('E51').Value = '  -2.83%  '
